# Bcrypt integration for secure password hashing
#
# The "Use bcrypt for hashing and storing passwords*" task (row 14) has been
# completed, so its related "In progress" status marker moves to "Done" and
# two more rows that previously had no status get marked "Done" as well.
#
# Net effect on the data:
#   - D9  ("In progress")      -> "Done"
#   - D15 (blank)               -> "Done"
#   - D20 (blank)               -> "Done"
# Once no cell references the shared string "In progress" any more, Excel
# drops it from the shared-strings table on save, which is why the saved
# workbook's string table shrinks/re-indexes automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "Done"
$ws.Range("D15").Value = "Done"
$ws.Range("D20").Value = "Done"

# Reflect the author's final cursor position/selection after making the edits.
$ws.Range("D21").Select() | Out-Null
